$d = $word.ActiveDocument

# The document currently has a single paragraph containing "1324567890".
# Add a new paragraph right after it containing "qwerty", inheriting the
# same run-level language formatting (en-US) via InsertParagraphAfter.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.Text = "qwerty"
